$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells AD1:AF1 - copy formatting (bold/border/centered) from an
# existing header cell (A1 uses style index 1) then overwrite the text.
$ws.Range("A1").Copy($ws.Range("AD1:AF1"))
$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Team record (Wins/Losses/Ties) repeated for every data row (2-51).
$lastRow = 51
for ($r = 2; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 30).Value = 101   # column AD
    $ws.Cells.Item($r, 31).Value = 61    # column AE
    $ws.Cells.Item($r, 32).Value = 1     # column AF
}
